$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the original frozen-header view (row 1 frozen, A2 top-left)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- Header text update ---
$ws.Range("C1").Value2 = "Effort Points (Hours)"

# --- Backlog grooming: re-ordered / re-estimated work items ---

# Row 2 - Project Proposal
$ws.Range("C2").Value2 = 13.0
$ws.Range("E2").Value2 = 45200.0
$ws.Range("F2").Value2 = "Past Due - ASAP"

# Row 3 - BRD (was "Coding Environment Setup")
$ws.Range("B3").Value2 = "BRD"
$ws.Range("C3").Value2 = 55.0
$ws.Range("D3").Value2 = "High "
$ws.Range("E3").Value2 = 45209.0
$ws.Range("F3").Value2 = 45210.0

# Row 4 - Team Review 1
$ws.Range("B4").Value2 = "Team Review 1"
$ws.Range("C4").Value2 = 21.0
$ws.Range("D4").Value2 = "Medium"
$ws.Range("E4").Value2 = 45209.0
$ws.Range("F4").Value2 = 45210.0

# Row 5 - Milestone 1
$ws.Range("B5").Value2 = "Milestone 1"
$ws.Range("C5").Value2 = 55.0
$ws.Range("D5").Value2 = "High "
$ws.Range("E5").Value2 = 45237.0
$ws.Range("F5").Value2 = 45238.0

# Row 6 - Team Review 2
$ws.Range("B6").Value2 = "Team Review 2"
$ws.Range("C6").Value2 = 21.0
$ws.Range("D6").Value2 = "Medium"
$ws.Range("E6").Value2 = 45237.0
$ws.Range("F6").Value2 = 45238.0

# Row 7 - Milestone 2
$ws.Range("B7").Value2 = "Milestone 2"
$ws.Range("C7").Value2 = 55.0
$ws.Range("D7").Value2 = "High "
$ws.Range("E7").Value2 = 45270.0
$ws.Range("F7").Value2 = 45271.0

# Row 8 - Team Review 3
$ws.Range("B8").Value2 = "Team Review 3"
$ws.Range("C8").Value2 = 21.0
$ws.Range("D8").Value2 = "Medium"
$ws.Range("E8").Value2 = 45270.0
$ws.Range("F8").Value2 = 45271.0

# Row 9 - Coding Environment Setup (was last "Past Due" row; now highlighted orange)
$ws.Range("B9").Value2 = "Coding Environment Setup"
$ws.Range("C9").Value2 = 8.0
$ws.Range("D9").Value2 = "Low "
$ws.Range("E9").Value2 = 45193.0
$ws.Range("F9").Value2 = 45194.0

# --- Formatting ---

# New right-alignment on the Effort figures for rows 2, 3 & 9
$ws.Range("C2").HorizontalAlignment = -4152
$ws.Range("C3").HorizontalAlignment = -4152
$ws.Range("C9").HorizontalAlignment = -4152

# Highlight the "Coding Environment Setup" row with the orange fill
$ws.Range("B9:F9").Interior.Color = 39423

# Restore the original selection
$ws.Range("B3").Select() | Out-Null
